$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.282.61'
$ws.Range('E2').Value = '  +1.05%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.918.70'
$ws.Range('E3').Value = '  +0.58%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.8100'
$ws.Range('E5').Value = '  +4.22%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '244.66'
$ws.Range('E6').Value = '  +1.32%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3246'
$ws.Range('E8').Value = '  +3.34%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '27.16'
$ws.Range('E9').Value = '  +4.86%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07097'

$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.7835'
$ws.Range('E11').Value = '  +6.24%  '

$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08095'
$ws.Range('E12').Value = '  +1.70%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.936.12'
$ws.Range('E13').Value = '  +1.54%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.382'
$ws.Range('E14').Value = '  +3.84%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '94.83'
$ws.Range('E15').Value = '  +2.56%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '30.280.78'
$ws.Range('E16').Value = '  +1.01%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.34'
$ws.Range('E17').Value = '  +3.65%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.027'
$ws.Range('E18').Value = '  +3.12%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '248.73'
$ws.Range('E19').Value = '  +1.80%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000007823'
$ws.Range('E20').Value = '  +1.61%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.175.32'
$ws.Range('E21').Value = '  +1.00%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.07%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.003'
$ws.Range('E23').Value = '  +0.04%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '7.737'
$ws.Range('E24').Value = '  +13.29%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.1631'
$ws.Range('E25').Value = '  +19.66%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.451'
$ws.Range('E26').Value = '  +2.15%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '167.77'
$ws.Range('E27').Value = '  -0.43%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.13'
$ws.Range('E28').Value = '  +1.44%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.123'
$ws.Range('E29').Value = '  +5.28%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.373'
$ws.Range('E30').Value = '  +0.22%  '

$ws.Range('E31').Value = '  +1.10%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.349'
$ws.Range('E32').Value = '  +1.10%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05661'
$ws.Range('E33').Value = '  +3.41%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.126'
$ws.Range('E34').Value = '  +1.46%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.303'
$ws.Range('E35').Value = '  +4.42%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7411'
$ws.Range('E36').Value = '  +1.63%  '

$ws.Range('E37').Value = '  +0.05%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.719'
$ws.Range('E38').Value = '  -0.45%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01952'
$ws.Range('E39').Value = '  +1.45%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.816'
$ws.Range('E40').Value = '  +1.12%  '

$ws.Range('E41').Value = '  +1.90%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '74.02'
$ws.Range('E42').Value = '  +3.45%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.984'
$ws.Range('E43').Value = '  -2.02%  '

$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.8538'
$ws.Range('E44').Value = '  +1.68%  '

$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.930'
$ws.Range('E45').Value = '  +3.72%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.048.56'
$ws.Range('E46').Value = '  +7.36%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.001'
$ws.Range('E47').Value = '  -0.11%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '103.00'
$ws.Range('E48').Value = '  +3.01%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.965'
$ws.Range('E49').Value = '  +2.23%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.622'
$ws.Range('E50').Value = '  +2.00%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.067.68'
$ws.Range('E51').Value = '  +0.59%  '
